$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing B2 value
$ws.Range("B2").Value = 24

# Add new rows 3-6 with column A (style copied from A2) and column B values
$ws.Range("A3").Value = 3
$ws.Range("B3").Value = 14

$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 9

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 5

$ws.Range("A6").Value = 2
$ws.Range("B6").Value = 5

# Copy the style from A2 to the new A3:A6 cells so they match (s="1")
$ws.Range("A2").Copy()
$ws.Range("A3:A6").PasteSpecial(-4122) # xlPasteFormats
